$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# Remove existing "总计" sheet; we'll recreate it after the new quarter sheet
# so sheetId ordering matches (2022-Q1 gets sheetId 3, 总计 gets sheetId 4).
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Delete()

# Create the new "2022-Q1" sheet right after "2021-Q2"
$q2 = $wb.Worksheets.Item("2021-Q2")
$q2.Range("B1:H2").Copy()
$newQ = $wb.Worksheets.Add($null, $q2)
$newQ.Name = "2022-Q1"
$newQ.Range("B1").PasteSpecial(-4122)

# Header row
$newQ.Range("B1").Value = "基金代码"
$newQ.Range("C1").Value = "基金名称"
$newQ.Range("D1").Value = "基金规模"
$newQ.Range("E1").Value = "股票总仓位"
$newQ.Range("F1").Value = "仓位占比"
$newQ.Range("G1").Value = "持有市值(亿元)"
$newQ.Range("H1").Value = "仓位排名"

# Data row 2 - text-typed numeric-looking values must keep "@" format to stay text
$q2.Range("A2").Copy()
$newQ.Range("A2").PasteSpecial(-4122)
$newQ.Range("A2").Value = 0
$newQ.Range("B2:G2").NumberFormat = "@"
$newQ.Range("B2").Value = "512780"
$newQ.Range("C2").Value = "广发中证京津冀协同发展主题ETF"
$newQ.Range("D2").Value = "0.13"
$newQ.Range("E2").Value = "98.52"
$newQ.Range("F2").Value = "3.30"
$newQ.Range("G2").Value = "0.0043"
$newQ.Range("H2").Value = 3

Write-Output "ok-2022Q1"

# Recreate "总计" sheet as the last sheet (after 2022-Q1), so it gets sheetId 4
$q2.Range("B1:D1").Copy()
$totalNew = $wb.Worksheets.Add($null, $newQ)
$totalNew.Name = "总计"
$totalNew.Range("B1").PasteSpecial(-4122)

$totalNew.Range("B1").Value = "日期"
$totalNew.Range("C1").Value = "持有数量(只)"
$totalNew.Range("D1").Value = "持有市值(亿元)"

$q2.Range("A2").Copy()
$totalNew.Range("A2:A4").PasteSpecial(-4122)

$totalNew.Range("A2").Value = 0
$totalNew.Range("B2").Value = "2022-Q1"
$totalNew.Range("C2").Value = 1
$totalNew.Range("D2").Value = 0

$totalNew.Range("A3").Value = 1
$totalNew.Range("B3").Value = "2021-Q2"
$totalNew.Range("C3").Value = 3
$totalNew.Range("D3").Value = 0.02

$totalNew.Range("A4").Value = 2
$totalNew.Range("B4").Value = "2020-Q4"
$totalNew.Range("C4").Value = 3
$totalNew.Range("D4").Value = 0.02

$wb.Worksheets.Item("2020-Q4").Activate()

Write-Output "ok-总计"
